$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1092.0526
$ws.Cells.Item(17, 10).Value = 1092.0526
$ws.Cells.Item(17, 12).Value = 3276.1578
$ws.Cells.Item(17, 14).Value = -3612.1578

$ws.Cells.Item(33, 8).Value = 725.25
$ws.Cells.Item(33, 9).Value = 183.04167
$ws.Cells.Item(33, 10).Value = 2351.875
$ws.Cells.Item(33, 11).Value = 183.04167
$ws.Cells.Item(33, 12).Value = 2351.875
$ws.Cells.Item(33, 13).Value = 45.95832999999999
$ws.Cells.Item(33, 14).Value = -2809.875

$ws.Cells.Item(64, 8).Value = 44523.543
$ws.Cells.Item(64, 9).Value = 85896.664
$ws.Cells.Item(64, 11).Value = 85896.664
$ws.Cells.Item(64, 13).Value = -85648.664

$ws.Cells.Item(67, 8).Value = 44523.543
$ws.Cells.Item(67, 9).Value = 85896.664
$ws.Cells.Item(67, 11).Value = 85896.664
$ws.Cells.Item(67, 13).Value = -85038.664

$ws.Cells.Item(107, 8).Value = 548.6923
$ws.Cells.Item(107, 9).Value = 582.35
$ws.Cells.Item(107, 10).Value = 436.5
$ws.Cells.Item(107, 11).Value = 582.35
$ws.Cells.Item(107, 12).Value = 436.5
$ws.Cells.Item(107, 13).Value = 1337.65
$ws.Cells.Item(107, 14).Value = -4276.5

$ws.Cells.Item(112, 8).Value = 1483.7142
$ws.Cells.Item(112, 9).Value = 490
$ws.Cells.Item(112, 10).Value = 1649.3334
$ws.Cells.Item(112, 11).Value = 1470
$ws.Cells.Item(112, 12).Value = 4948.0002
$ws.Cells.Item(112, 13).Value = -362
$ws.Cells.Item(112, 14).Value = -7164.0002

$ws.Cells.Item(125, 8).Value = 4262.125
$ws.Cells.Item(125, 10).Value = 3437.8572
$ws.Cells.Item(125, 12).Value = 30940.7148
$ws.Cells.Item(125, 14).Value = -35860.7148

$ws.Cells.Item(138, 8).Value = 1625.683
$ws.Cells.Item(138, 9).Value = 1789.8667
$ws.Cells.Item(138, 10).Value = 1530.9615
$ws.Cells.Item(138, 11).Value = 5369.6001
$ws.Cells.Item(138, 12).Value = 4592.8845
$ws.Cells.Item(138, 13).Value = -229.6000999999997
$ws.Cells.Item(138, 14).Value = -14872.8845

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 51399.35
$ws.Cells.Item(2, 9).Value = 1467.125
$ws.Cells.Item(2, 11).Value = 1467.125
$ws.Cells.Item(2, 13).Value = -1354.125

$ws.Cells.Item(32, 8).Value = 28067.455
$ws.Cells.Item(32, 9).Value = 6055.828
$ws.Cells.Item(32, 10).Value = 380253.5
$ws.Cells.Item(32, 11).Value = 6055.828
$ws.Cells.Item(32, 12).Value = 380253.5
$ws.Cells.Item(32, 13).Value = -5768.828
$ws.Cells.Item(32, 14).Value = -380827.5

$ws.Cells.Item(97, 8).Value = 25859.85
$ws.Cells.Item(97, 9).Value = 32840.71
$ws.Cells.Item(97, 10).Value = 1814.6666
$ws.Cells.Item(97, 11).Value = 32840.71
$ws.Cells.Item(97, 12).Value = 1814.6666
$ws.Cells.Item(97, 13).Value = -32344.71
$ws.Cells.Item(97, 14).Value = -2806.6666

$ws.Cells.Item(102, 8).Value = 73362.36
$ws.Cells.Item(102, 9).Value = 201898
$ws.Cells.Item(102, 10).Value = 1953.6666
$ws.Cells.Item(102, 11).Value = 201898
$ws.Cells.Item(102, 12).Value = 1953.6666
$ws.Cells.Item(102, 13).Value = -200276
$ws.Cells.Item(102, 14).Value = -5197.6666

$ws.Cells.Item(110, 8).Value = 37078570
$ws.Cells.Item(110, 9).Value = 47672188
$ws.Cells.Item(110, 11).Value = 47672188
$ws.Cells.Item(110, 13).Value = -47670143

$ws.Cells.Item(116, 8).Value = 51399.35
$ws.Cells.Item(116, 9).Value = 1467.125
$ws.Cells.Item(116, 11).Value = 1467.125
$ws.Cells.Item(116, 13).Value = 826.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 51399.35
$ws.Cells.Item(3, 9).Value = 1467.125
$ws.Cells.Item(3, 11).Value = 1467.125
$ws.Cells.Item(3, 13).Value = -1353.125

$ws.Cells.Item(20, 8).Value = 127887.5
$ws.Cells.Item(20, 9).Value = 145728.58
$ws.Cells.Item(20, 10).Value = 3000
$ws.Cells.Item(20, 11).Value = 145728.58
$ws.Cells.Item(20, 12).Value = 3000
$ws.Cells.Item(20, 13).Value = -145481.58
$ws.Cells.Item(20, 14).Value = -3494

$ws.Cells.Item(80, 8).Value = 1134.8214
$ws.Cells.Item(80, 10).Value = 1213.75
$ws.Cells.Item(80, 12).Value = 1213.75
$ws.Cells.Item(80, 14).Value = -3209.75

$ws.Cells.Item(83, 8).Value = 1134.8214
$ws.Cells.Item(83, 10).Value = 1213.75
$ws.Cells.Item(83, 12).Value = 6068.75
$ws.Cells.Item(83, 14).Value = -16052.75

$ws.Cells.Item(86, 8).Value = 162143.42
$ws.Cells.Item(86, 9).Value = 188667.33
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 188667.33
$ws.Cells.Item(86, 12).Value = 3000
$ws.Cells.Item(86, 13).Value = -187544.33
$ws.Cells.Item(86, 14).Value = -5246

$ws.Cells.Item(89, 8).Value = 162143.42
$ws.Cells.Item(89, 9).Value = 188667.33
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 943336.6499999999
$ws.Cells.Item(89, 12).Value = 15000
$ws.Cells.Item(89, 13).Value = -937720.6499999999
$ws.Cells.Item(89, 14).Value = -26232

$ws.Cells.Item(94, 8).Value = 604
$ws.Cells.Item(94, 9).Value = 462.33334
$ws.Cells.Item(94, 10).Value = 944
$ws.Cells.Item(94, 11).Value = 462.33334
$ws.Cells.Item(94, 12).Value = 944
$ws.Cells.Item(94, 13).Value = -11.33334000000002
$ws.Cells.Item(94, 14).Value = -1846

$ws.Cells.Item(107, 8).Value = 71462160
$ws.Cells.Item(107, 9).Value = 100046150
$ws.Cells.Item(107, 10).Value = 2188.5
$ws.Cells.Item(107, 11).Value = 100046150
$ws.Cells.Item(107, 12).Value = 2188.5
$ws.Cells.Item(107, 13).Value = -100044230
$ws.Cells.Item(107, 14).Value = -6028.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2195.0344
$ws.Cells.Item(58, 9).Value = 1668.1818
$ws.Cells.Item(58, 10).Value = 2517
$ws.Cells.Item(58, 11).Value = 1668.1818
$ws.Cells.Item(58, 12).Value = 2517
$ws.Cells.Item(58, 13).Value = -1465.1818
$ws.Cells.Item(58, 14).Value = -2923

$ws.Cells.Item(94, 8).Value = 1044.2916
$ws.Cells.Item(94, 9).Value = 1065.3334
$ws.Cells.Item(94, 10).Value = 1037.2778
$ws.Cells.Item(94, 11).Value = 1065.3334
$ws.Cells.Item(94, 12).Value = 1037.2778
$ws.Cells.Item(94, 13).Value = -614.3334
$ws.Cells.Item(94, 14).Value = -1939.2778

$ws.Cells.Item(99, 8).Value = 13141.8
$ws.Cells.Item(99, 9).Value = 3803.4285
$ws.Cells.Item(99, 10).Value = 34931.332
$ws.Cells.Item(99, 11).Value = 3803.4285
$ws.Cells.Item(99, 12).Value = 34931.332
$ws.Cells.Item(99, 13).Value = -2305.4285
$ws.Cells.Item(99, 14).Value = -37927.332

$ws.Cells.Item(105, 8).Value = 1184.3889
$ws.Cells.Item(105, 9).Value = 1136.5834
$ws.Cells.Item(105, 10).Value = 1280
$ws.Cells.Item(105, 11).Value = 1136.5834
$ws.Cells.Item(105, 12).Value = 1280
$ws.Cells.Item(105, 13).Value = 610.4166
$ws.Cells.Item(105, 14).Value = -4774

$ws.Cells.Item(126, 8).Value = 13141.8
$ws.Cells.Item(126, 9).Value = 3803.4285
$ws.Cells.Item(126, 10).Value = 34931.332
$ws.Cells.Item(126, 11).Value = 11410.2855
$ws.Cells.Item(126, 12).Value = 104793.996
$ws.Cells.Item(126, 13).Value = -8940.2855
$ws.Cells.Item(126, 14).Value = -109733.996

$ws.Cells.Item(131, 8).Value = 19075.191
$ws.Cells.Item(131, 10).Value = 19075.191
$ws.Cells.Item(131, 12).Value = 19075.191
$ws.Cells.Item(131, 14).Value = -29155.191

$ws.Cells.Item(136, 8).Value = 2195.0344
$ws.Cells.Item(136, 9).Value = 1668.1818
$ws.Cells.Item(136, 10).Value = 2517
$ws.Cells.Item(136, 11).Value = 5004.5454
$ws.Cells.Item(136, 12).Value = 7551
$ws.Cells.Item(136, 13).Value = -2454.5454
$ws.Cells.Item(136, 14).Value = -12651

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1622.6666
$ws.Cells.Item(4, 9).Value = 146.66667
$ws.Cells.Item(4, 10).Value = 1991.6666
$ws.Cells.Item(4, 11).Value = 440.00001
$ws.Cells.Item(4, 12).Value = 5974.9998
$ws.Cells.Item(4, 13).Value = -328.00001
$ws.Cells.Item(4, 14).Value = -6198.9998

$ws.Cells.Item(34, 8).Value = 1142.1428
$ws.Cells.Item(34, 10).Value = 1299.1666
$ws.Cells.Item(34, 12).Value = 3897.4998
$ws.Cells.Item(34, 14).Value = -4065.4998

$ws.Cells.Item(107, 8).Value = 480048.7
$ws.Cells.Item(107, 9).Value = 800
$ws.Cells.Item(107, 10).Value = 927347.4399999999
$ws.Cells.Item(107, 11).Value = 2400
$ws.Cells.Item(107, 12).Value = 2782042.32
$ws.Cells.Item(107, 13).Value = -480
$ws.Cells.Item(107, 14).Value = -2785882.32

$ws.Cells.Item(121, 8).Value = 6572.04
$ws.Cells.Item(121, 9).Value = 15509.667
$ws.Cells.Item(121, 10).Value = 5353.273
$ws.Cells.Item(121, 11).Value = 46529.001
$ws.Cells.Item(121, 12).Value = 16059.819
$ws.Cells.Item(121, 13).Value = -45219.001
$ws.Cells.Item(121, 14).Value = -18679.819

$ws.Cells.Item(131, 8).Value = 832.5
$ws.Cells.Item(131, 10).Value = 843.4693600000001
$ws.Cells.Item(131, 12).Value = 2530.40808
$ws.Cells.Item(131, 14).Value = -12610.40808

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(45, 8).Value = 180250
$ws.Cells.Item(45, 10).Value = 180250
$ws.Cells.Item(45, 12).Value = 180250
$ws.Cells.Item(45, 14).Value = -181368

$ws.Cells.Item(51, 8).Value = 113500
$ws.Cells.Item(51, 10).Value = 113500
$ws.Cells.Item(51, 12).Value = 113500
$ws.Cells.Item(51, 14).Value = -114518

$ws.Cells.Item(70, 8).Value = 55276.3
$ws.Cells.Item(70, 10).Value = 5710.5
$ws.Cells.Item(70, 12).Value = 5710.5
$ws.Cells.Item(70, 14).Value = -6250.5

$ws.Cells.Item(73, 8).Value = 55276.3
$ws.Cells.Item(73, 10).Value = 5710.5
$ws.Cells.Item(73, 12).Value = 5710.5
$ws.Cells.Item(73, 14).Value = -7582.5

$ws.Cells.Item(97, 8).Value = 142860240
$ws.Cells.Item(97, 9).Value = 166669950
$ws.Cells.Item(97, 10).Value = 2000
$ws.Cells.Item(97, 11).Value = 166669950
$ws.Cells.Item(97, 12).Value = 2000
$ws.Cells.Item(97, 13).Value = -166669454
$ws.Cells.Item(97, 14).Value = -2992

$ws.Cells.Item(132, 8).Value = 2300.1843
$ws.Cells.Item(132, 9).Value = 2137.3125
$ws.Cells.Item(132, 10).Value = 3168.8333
$ws.Cells.Item(132, 11).Value = 6411.9375
$ws.Cells.Item(132, 12).Value = 9506.499899999999
$ws.Cells.Item(132, 13).Value = -3881.9375
$ws.Cells.Item(132, 14).Value = -14566.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 920582.25
$ws.Cells.Item(46, 9).Value = 437.8
$ws.Cells.Item(46, 10).Value = 1687369.4
$ws.Cells.Item(46, 11).Value = 437.8
$ws.Cells.Item(46, 12).Value = 1687369.4
$ws.Cells.Item(46, 13).Value = -249.8
$ws.Cells.Item(46, 14).Value = -1687745.4

$ws.Cells.Item(69, 8).Value = 37000
$ws.Cells.Item(69, 10).Value = 37000
$ws.Cells.Item(69, 12).Value = 37000
$ws.Cells.Item(69, 14).Value = -38622

$ws.Cells.Item(72, 8).Value = 37000
$ws.Cells.Item(72, 10).Value = 37000
$ws.Cells.Item(72, 12).Value = 111000
$ws.Cells.Item(72, 14).Value = -119112

$ws.Cells.Item(136, 8).Value = 1845.75
$ws.Cells.Item(136, 9).Value = 1707.24
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 11).Value = 5121.72
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = -2571.72
$ws.Cells.Item(136, 14).Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 710.7917
$ws.Cells.Item(113, 9).Value = 538.26666
$ws.Cells.Item(113, 10).Value = 998.3333
$ws.Cells.Item(113, 11).Value = 1614.79998
$ws.Cells.Item(113, 12).Value = 2994.9999
$ws.Cells.Item(113, 13).Value = 555.20002
$ws.Cells.Item(113, 14).Value = -7334.9999
